$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 356, pushing the existing data
# (old rows 356..451) down to rows 358..453.
$ws.Rows.Item(356).Insert()
$ws.Rows.Item(356).Insert()

# Populate the two newly inserted rows (356 and 357) with new price entries.
# Columns A,B,C,E,F,G,H,I,J share the same boilerplate values across every
# row in this sheet (Mercado/Region/Producto metadata).

# Row 356: Florida King / Especial
$ws.Cells.Item(356, 1).Value2 = 10
$ws.Cells.Item(356, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(356, 3).Value = "La Araucanía"
$ws.Cells.Item(356, 4).Value2 = 45244
$ws.Cells.Item(356, 5).Value2 = 9
$ws.Cells.Item(356, 6).Value = "Fruta"
$ws.Cells.Item(356, 7).Value2 = 100103
$ws.Cells.Item(356, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(356, 9).Value2 = 100103004
$ws.Cells.Item(356, 10).Value = "Durazno"
$ws.Cells.Item(356, 11).Value = "Florida King"
$ws.Cells.Item(356, 12).Value = "Especial"
$ws.Cells.Item(356, 13).Value2 = 40
$ws.Cells.Item(356, 14).Value2 = 20000
$ws.Cells.Item(356, 15).Value2 = 20000
$ws.Cells.Item(356, 16).Value2 = 20000
$ws.Cells.Item(356, 17).Value = "$/bandeja 10 kilos granel"
$ws.Cells.Item(356, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(356, 19).Value2 = 2000
$ws.Cells.Item(356, 20).Value2 = 10

# Row 357: Florida King / Primera
$ws.Cells.Item(357, 1).Value2 = 10
$ws.Cells.Item(357, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(357, 3).Value = "La Araucanía"
$ws.Cells.Item(357, 4).Value2 = 45244
$ws.Cells.Item(357, 5).Value2 = 9
$ws.Cells.Item(357, 6).Value = "Fruta"
$ws.Cells.Item(357, 7).Value2 = 100103
$ws.Cells.Item(357, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(357, 9).Value2 = 100103004
$ws.Cells.Item(357, 10).Value = "Durazno"
$ws.Cells.Item(357, 11).Value = "Florida King"
$ws.Cells.Item(357, 12).Value = "Primera"
$ws.Cells.Item(357, 13).Value2 = 100
$ws.Cells.Item(357, 14).Value2 = 16000
$ws.Cells.Item(357, 15).Value2 = 16000
$ws.Cells.Item(357, 16).Value2 = 16000
$ws.Cells.Item(357, 17).Value = "$/bandeja 10 kilos granel"
$ws.Cells.Item(357, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(357, 19).Value2 = 1600
$ws.Cells.Item(357, 20).Value2 = 10
